$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.093.25"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "1.780.50"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("E4").Value = "  +0.18%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "225.35"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("E7").Value = "  +0.18%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "31.79"
$ws.Range("E8").Value = "  -1.20%  "

$ws.Range("E9").Value = "  -1.43%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0685"
$ws.Range("E10").Value = "  -0.01%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0947"
$ws.Range("E11").Value = "  +0.70%  "

$ws.Range("D12").Value = "2.037.10"
$ws.Range("E12").Value = "  -0.46%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "10.92"
$ws.Range("E13").Value = "  -3.37%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.795.31"
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("D15").Value = "34.071.85"
$ws.Range("E15").Value = "  +0.07%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.620"
$ws.Range("E16").Value = "  -0.54%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "4.18"
$ws.Range("E17").Value = "  -0.29%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "67.51"
$ws.Range("E18").Value = "  -0.59%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "245.29"
$ws.Range("E19").Value = "  +1.28%  "

$ws.Range("D20").Value = "0.0₃0786"
$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("E21").Value = "  +0.31%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.85"
$ws.Range("E22").Value = "  +1.16%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.09"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("E24").Value = "  -1.38%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "161.66"
$ws.Range("E25").Value = "  -0.18%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.09"
$ws.Range("E26").Value = "  -0.68%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "16.22"
$ws.Range("E27").Value = "  +0.24%  "

$ws.Range("E28").Value = "  +0.37%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Range("E29").Value = "  +0.22%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.22"
$ws.Range("E30").Value = "  -1.01%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0516"
$ws.Range("E31").Value = "  +0.07%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.70"
$ws.Range("E32").Value = "  +1.37%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.71"
$ws.Range("E33").Value = "  +2.56%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.79"
$ws.Range("E34").Value = "  -2.32%  "

$ws.Range("D35").Value = "1.444.42"
$ws.Range("E35").Value = "  +2.94%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.45"
$ws.Range("E36").Value = "  +4.91%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.648"
$ws.Range("E37").Value = "  -0.45%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0190"
$ws.Range("E38").Value = "  +0.86%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.04"
$ws.Range("E39").Value = "  -0.46%  "

$ws.Range("E40").Value = "  +1.61%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "80.56"
$ws.Range("E41").Value = "  +0.96%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "2.72"
$ws.Range("E42").Value = "  +1.29%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.914"
$ws.Range("E43").Value = "  -0.37%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "13.70"
$ws.Range("E44").Value = "  +0.62%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0517"
$ws.Range("E45").Value = "  +2.09%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "6.05"
$ws.Range("E46").Value = "  -1.37%  "

$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.08"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0133"
$ws.Range("E48").Value = "  -4.49%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.937.20"
$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "104.23"
$ws.Range("E50").Value = "  -3.05%  "

$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Range("E51").Value = "  +0.19%  "

